$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ('K') values per row, regenerated per the author's commit
# ("regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals").
$kValues = @{
    2 = 2
    3 = 3
    4 = 1
    5 = 0
    6 = 4
    7 = 1
    8 = 0
    9 = 3
    10 = 0
    11 = 4
    12 = 3
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 4
    27 = 1
    28 = 0
    29 = 8
    30 = 0
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 2
    38 = 4
    39 = 2
    40 = 3
    41 = 4
    42 = 0
    43 = 2
    44 = 2
    45 = 2
    46 = 1
    47 = 2
    48 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
